$wb = $excel.ActiveWorkbook

# Update zh-cn sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 03:17:14"
$wsZhCn.Range("H2").Value = "2016-03-18 03:17:53"

# Update de-de sheet: Correspond Handoff Datetime (E2) and Correspond Handback DateTime (H2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 03:17:22"
$wsDeDe.Range("H2").Value = "2016-03-18 03:18:07"
